# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap Eslovaquia / Croacia rows (data follows the country, rows keep their
#     original position in the sheet but the country label + stats move) ---
$ws.Range("A88").Value = "Croacia"
$ws.Range("A89").Value = "Eslovaquia"

# --- Swap Sudan / Georgia rows ---
$ws.Range("A100").Value = "Georgia"
$ws.Range("A101").Value = "Sudan"

# --- Update "last updated" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 15 de Octubre de 2020 a las 10:18"

# --- Row 4 (Estados Unidos) ---
$ws.Range("B4").Value = 8150383
$ws.Range("C4").Value = 340
$ws.Range("D4").Value = 5279651
$ws.Range("E4").Value = 2648882
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 7
$ws.Range("H4").Value = 221850

# --- Row 7 (Rusia) ---
$ws.Range("B7").Value = 1354163
$ws.Range("C7").Value = 13754
$ws.Range("D7").Value = 1048097
$ws.Range("E7").Value = 282575
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 286
$ws.Range("H7").Value = 23491

# --- Row 21 (Filipinas) ---
$ws.Range("B21").Value = 348698
$ws.Range("C21").Value = 2261
$ws.Range("D21").Value = 294161
$ws.Range("E21").Value = 48040
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 50
$ws.Range("H21").Value = 6497

# --- Row 35 (Chequia) ---
$ws.Range("D35").Value = 85588
$ws.Range("E35").Value = 52999

# --- Row 64 (Singapur) ---
$ws.Range("B64").Value = 57892
$ws.Range("C64").Value = 3
$ws.Range("E64").Value = 112

# --- Row 88 (now Croacia) ---
$ws.Range("B88").Value = 22534
$ws.Range("C88").Value = 793
$ws.Range("D88").Value = 18628
$ws.Range("E88").Value = 3562
$ws.Range("G88").Value = 10
$ws.Range("H88").Value = 344

# --- Row 89 (now Eslovaquia) ---
$ws.Range("B89").Value = 22296
$ws.Range("D89").Value = 6709
$ws.Range("E89").Value = 15521
$ws.Range("H89").Value = 66

# --- Row 100 (now Georgia) ---
$ws.Range("B100").Value = 14440
$ws.Range("C100").Value = 919
$ws.Range("D100").Value = 7367
$ws.Range("E100").Value = 6960
$ws.Range("G100").Value = 4
$ws.Range("H100").Value = 113

# --- Row 101 (now Sudan) ---
$ws.Range("B101").Value = 13691
$ws.Range("D101").Value = 6764
$ws.Range("E101").Value = 6091
$ws.Range("H101").Value = 836

# --- Row 121 (Mauritania) ---
$ws.Range("B121").Value = 6760
$ws.Range("C121").Value = 255
$ws.Range("D121").Value = 2983
$ws.Range("E121").Value = 3667
$ws.Range("G121").Value = 1
$ws.Range("H121").Value = 110

# --- Row 141 (Reunion) ---
$ws.Range("B141").Value = 3980
$ws.Range("C141").Value = 33
$ws.Range("D141").Value = 3093

$wb.Save()
